$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 (species observation records) get cyclically shifted up by one:
# new row2 = old row3, new row3 = old row4, new row4 = old row2.
# Only columns A, B, E, F, G, H, Q, R, Z, AB actually differ between the rows;
# the rest are identical across rows 2-4 so they are left untouched.

$cols = @("A","B","E","F","G","H","Q","R","Z","AB")

$old2 = @{}
$old3 = @{}
$old4 = @{}
foreach ($col in $cols) {
    $old2[$col] = $ws.Range($col + "2").Value()
    $old3[$col] = $ws.Range($col + "3").Value()
    $old4[$col] = $ws.Range($col + "4").Value()
}

foreach ($col in $cols) {
    $ws.Range($col + "2").Value = $old3[$col]
    $ws.Range($col + "3").Value = $old4[$col]
    $ws.Range($col + "4").Value = $old2[$col]
}
